$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A93").NumberFormat = "@"
$ws.Range("A93").Value = "2021/10/13"
$ws.Range("A93").NumberFormat = "yyyy/mm/dd"
$ws.Range("B93").Value = 83.2
$ws.Range("C93").Value = 83.2
$ws.Range("D93").Value = 0.97
$ws.Range("E93").Value = 0.98

$ws.Range("A94").Select() | Out-Null
